$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values per row: B = Year1, C = Year2, D = Total
$data = @{
    2  = @(-38293.06, -31954.07, -70247.13)
    3  = @(-38819.41, -38946.09, -77765.5)
    4  = @(-15438.74, -15008.15, -30446.89)
    5  = @(-31167.05, -37792.51, -68959.56)
    6  = @(-16132.08, -18620.19, -34752.27)
    7  = @(-54028.1,  -73444.08, -127472.18)
    8  = @(-27644.43, -29930.08, -57574.51)
    9  = @(-15629.56, -19907.67, -35537.23)
    10 = @(-237152.43, -265602.84, -502755.27)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
}
